$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the year header row (row 4): was 2007-2019 across D:P, becomes 2015-2021 across D:J ---
$ws.Range("D4").Value = 2015
$ws.Range("E4").Value = 2016
$ws.Range("F4").Value = 2017
$ws.Range("G4").Value = 2018
$ws.Range("H4").Value = 2019
$ws.Range("I4").Value = 2020
$ws.Range("J4").Value = 2021

# --- Row 5 data: new values; D5:G5 need to switch from the old one-off style to the
#     same style already used by H5:J5 (copy number format/font/etc from H5) ---
$ws.Range("H5").Copy()
$ws.Range("D5:G5").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("D5").Value = 2.2197193775563164
$ws.Range("E5").Value = 2.1235271668715399
$ws.Range("F5").Value = 2.7818537161298167
$ws.Range("G5").Value = 6.7272960584548969
$ws.Range("H5").Value = 5.1525830614767187
$ws.Range("I5").Value = 4.4774536255935971
$ws.Range("J5").Value = 4.6024666695867751

# --- Row 6 data: new values; E6:J6 need to switch to the style already used further right
#     in the old layout (M6's style) while D6 keeps its own style ---
$ws.Range("M6").Copy()
$ws.Range("E6:J6").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("D6").Value = 2.2322863217945752
$ws.Range("E6").Value = 2.8603553109638966
$ws.Range("F6").Value = 3.113207036164539
$ws.Range("G6").Value = 6.2970593463100784
$ws.Range("H6").Value = 4.8617746111834492
$ws.Range("I6").Value = 2.6715092780025032
$ws.Range("J6").Value = 4.3694509108608912

# --- The table no longer carries the 2007-2014 / 2018-2019 years (columns K:P) at all ---
$ws.Range("K1:P6").Clear()

# --- Column widths for D:J change from mixed/default widths to a uniform custom width ---
$ws.Columns("D:J").ColumnWidth = 8.65

# --- Selection moves ---
$ws.Range("K16").Select() | Out-Null
